$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 415.08334
$ws.Range("I58").Value = 298.27274
$ws.Range("J58").Value = 1700
$ws.Range("K58").Value = 894.81822
$ws.Range("L58").Value = 5100
$ws.Range("M58").Value = -744.81822
$ws.Range("N58").Value = -5400

$ws.Range("H62").Value = 41115.195
$ws.Range("I62").Value = 57128.55
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 57128.55
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -56504.55
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 41115.195
$ws.Range("I65").Value = 57128.55
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 285642.75
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -282522.75
$ws.Range("N65").Value = -66240

$ws.Range("H112").Value = 55836.316
$ws.Range("J112").Value = 58917.5
$ws.Range("L112").Value = 176752.5
$ws.Range("N112").Value = -178968.5

$ws.Range("H127").Value = 1207
$ws.Range("J127").Value = 1700
$ws.Range("L127").Value = 5100
$ws.Range("N127").Value = -15020

$ws.Range("H132").Value = 20897700
$ws.Range("I132").Value = 23811450
$ws.Range("K132").Value = 71434350
$ws.Range("M132").Value = -71431820

$ws.Range("H138").Value = 3364.9375
$ws.Range("J138").Value = 5690.727
$ws.Range("L138").Value = 17072.181
$ws.Range("N138").Value = -27352.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1985.125
$ws.Range("I5").Value = 2155.5
$ws.Range("K5").Value = 2155.5
$ws.Range("M5").Value = -2043.5

$ws.Range("H61").Value = 5840.8184
$ws.Range("I61").Value = 5713.857
$ws.Range("J61").Value = 6063
$ws.Range("K61").Value = 5713.857
$ws.Range("L61").Value = 6063
$ws.Range("M61").Value = -5501.857
$ws.Range("N61").Value = -6487

$ws.Range("H97").Value = 2057.2354
$ws.Range("I97").Value = 2057.2354
$ws.Range("K97").Value = 2057.2354
$ws.Range("M97").Value = -1561.2354

$ws.Range("H132").Value = 4738.4287
$ws.Range("I132").Value = 4250
$ws.Range("K132").Value = 12750
$ws.Range("M132").Value = -10220

$ws.Range("H136").Value = 5840.8184
$ws.Range("I136").Value = 5713.857
$ws.Range("J136").Value = 6063
$ws.Range("K136").Value = 17141.571
$ws.Range("L136").Value = 18189
$ws.Range("M136").Value = -14591.571
$ws.Range("N136").Value = -23289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1985.125
$ws.Range("I4").Value = 2155.5
$ws.Range("K4").Value = 2155.5
$ws.Range("M4").Value = -2040.5

$ws.Range("H22").Value = 412957.25
$ws.Range("I22").Value = 627.6667
$ws.Range("K22").Value = 627.6667
$ws.Range("M22").Value = -454.6667

$ws.Range("H94").Value = 1557.5714
$ws.Range("I94").Value = 1418.8572
$ws.Range("J94").Value = 2112.4285
$ws.Range("K94").Value = 1418.8572
$ws.Range("L94").Value = 2112.4285
$ws.Range("M94").Value = -967.8571999999999
$ws.Range("N94").Value = -3014.4285

$ws.Range("H99").Value = 2921.5173
$ws.Range("I99").Value = 1522.091
$ws.Range("J99").Value = 7319.7144
$ws.Range("K99").Value = 1522.091
$ws.Range("L99").Value = 7319.7144
$ws.Range("M99").Value = -24.09099999999989
$ws.Range("N99").Value = -10315.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 64421.293
$ws.Range("I31").Value = 113466.445
$ws.Range("J31").Value = 9245.5
$ws.Range("K31").Value = 113466.445
$ws.Range("L31").Value = 9245.5
$ws.Range("M31").Value = -113171.445
$ws.Range("N31").Value = -9835.5

$ws.Range("H34").Value = 64421.293
$ws.Range("I34").Value = 113466.445
$ws.Range("J34").Value = 9245.5
$ws.Range("K34").Value = 113466.445
$ws.Range("L34").Value = 9245.5
$ws.Range("M34").Value = -113264.445
$ws.Range("N34").Value = -9649.5

$ws.Range("H132").Value = 3879.7856
$ws.Range("I132").Value = 3778.2307
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 11334.6921
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -8804.6921
$ws.Range("N132").Value = -20660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 10044.667
$ws.Range("J104").Value = 10044.667
$ws.Range("L104").Value = 30134.001
$ws.Range("N104").Value = -35376.001

$ws.Range("H124").Value = 23213.846
$ws.Range("I124").Value = 7277.8335
$ws.Range("J124").Value = 27994.65
$ws.Range("K124").Value = 21833.5005
$ws.Range("L124").Value = 83983.95000000001
$ws.Range("M124").Value = -16923.5005
$ws.Range("N124").Value = -93803.95000000001

$ws.Range("H132").Value = 1612.5
$ws.Range("I132").Value = 980
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 8820
$ws.Range("L132").Value = 24000.0003
$ws.Range("M132").Value = -6290
$ws.Range("N132").Value = -29060.0003

$ws.Range("H138").Value = 16675644
$ws.Range("J138").Value = 8055.3184
$ws.Range("L138").Value = 24165.9552
$ws.Range("N138").Value = -34445.9552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 3187.5
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 3916.6667
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 3916.6667
$ws.Range("M19").Value = -712
$ws.Range("N19").Value = -4492.6667

$ws.Range("H102").Value = 3164.3333
$ws.Range("I102").Value = 3197.2
$ws.Range("K102").Value = 3197.2
$ws.Range("M102").Value = -1575.2

$ws.Range("H132").Value = 3183
$ws.Range("J132").Value = 7000
$ws.Range("L132").Value = 21000
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 1000000
$ws.Range("I23").Value = 1000000
$ws.Range("K23").Value = 1000000
$ws.Range("M23").Value = -999770

$ws.Range("H40").Value = 7801
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864

$ws.Range("H58").Value = 3325
$ws.Range("I58").Value = 575
$ws.Range("K58").Value = 575
$ws.Range("M58").Value = -315

$ws.Range("H132").Value = 2628.2646
$ws.Range("I132").Value = 2230.0312
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 6690.0936
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -4160.0936
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10562.1875
$ws.Range("I81").Value = 27675.25
$ws.Range("J81").Value = 4857.8335
$ws.Range("K81").Value = 55350.5
$ws.Range("L81").Value = 9715.666999999999
$ws.Range("M81").Value = -54289.5
$ws.Range("N81").Value = -11837.667

$ws.Range("H84").Value = 10562.1875
$ws.Range("I84").Value = 27675.25
$ws.Range("J84").Value = 4857.8335
$ws.Range("K84").Value = 276752.5
$ws.Range("L84").Value = 48578.335
$ws.Range("M84").Value = -271448.5
$ws.Range("N84").Value = -59186.335

$ws.Range("H107").Value = 539.44446
$ws.Range("I107").Value = 544.25
$ws.Range("K107").Value = 1632.75
$ws.Range("M107").Value = 287.25

$ws.Range("H126").Value = 3443.7778
$ws.Range("I126").Value = 2713.4285
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 8140.2855
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -5670.2855
$ws.Range("N126").Value = -22940

$ws.Range("H132").Value = 4515.6787
$ws.Range("I132").Value = 4151.7827
$ws.Range("J132").Value = 6189.6
$ws.Range("K132").Value = 12455.3481
$ws.Range("L132").Value = 18568.8
$ws.Range("M132").Value = -9925.348099999999
$ws.Range("N132").Value = -23628.8

